$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark PBI #10 (row 6) as Done - the Effort-Forecast feature is now finished.
$ws.Range("E6").Value = "Done"

# Insert two new rows for the newly created backlog items (PBI 11 and PBI 12),
# pushing the existing rows 7-11 down to rows 9-13.
$ws.Rows("7:8").Insert()

# Row 7: PBI 11 - accumulated Effort-Forecast per sprint
$ws.Range("A7").Value = 11
$ws.Range("B7").Value = "Als PO möchte ich den akkumulieren Effort-Forecast bis zu jedem Sprint sehen können."
$ws.Range("C7").Value = "Akzeptanzkriterien:`n- Die Geschwindigkeiten werden wie bei [10] berechnet.`n- Der akkumulierte Forecast ergibt sich aus dem akkumulierten Forecast des vorherigen Sprints plus dem`nEffort-Forecast dieses Sprints"
$ws.Range("C7").WrapText = $true
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = "Todo"
$ws.Range("F7").Value = "Sprint 3"
$ws.Rows("7").RowHeight = 120

# Row 8: PBI 12 - Sprint-Forecast per PBI
$ws.Range("A8").Value = 12
$ws.Range("B8").Value = "Als PO möchte ich für alle PBIs einen Forecast bis zu welchem Sprint dieses fertiggestellt wird."
$ws.Range("C8").Value = "Akzeptanzkritrien:`n- Hierzu werden die akkumulierten Effort-Forecast aus [11] verwendet.`n- Der Sprint-Forecast wird im PBL dargestellt"
$ws.Range("C8").WrapText = $true
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = "Todo"
$ws.Range("F8").Value = "Sprint 3"
$ws.Rows("8").RowHeight = 90

# Move the visible selection to reflect the newly added rows.
$ws.Range("F9").Select() | Out-Null
